# Localization sheet update: add several new UI-layout / art related
# language keys (loading/options/organism-editor strings) and rename the
# "categoryBody" value from "Shape" to "Body".
#
# The statements below are intentionally ordered to match the order the
# strings were actually typed into the workbook (new rows were inserted in
# batches, filled out of final top-to-bottom order), which is what drives
# the order entries are appended to xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shift the existing "body / essential / motility / metabolism" key
#     table rows down (they keep their original text, just a new row
#     number), making room for the freshly inserted rows above them.

$ws.Range("A25").Value = "testBodyCapsule"
$ws.Range("B25").Value = "Capsule"
$ws.Range("A26").Value = "testBodySphere"
$ws.Range("B26").Value = "Sphere"
$ws.Range("A27").Value = "categoryBody"
$ws.Range("A28").Value = "categoryCellStructure"
$ws.Range("B28").Value = "Structure"
$ws.Range("A29").Value = "categoryMotility"
$ws.Range("B29").Value = "Motility"
$ws.Range("A30").Value = "categoryMetabolism"
$ws.Range("B30").Value = "Metabolism"
$ws.Range("A31").Value = "essentialNucleoid"
$ws.Range("B31").Value = "Nucleoid"
$ws.Range("A32").Value = "essentialRibosome"
$ws.Range("B32").Value = "Ribosome"
$ws.Range("A33").Value = "essentialPlasmid"
$ws.Range("B33").Value = "Plasmid DNA"
$ws.Range("A34").Value = "bodyBacillus"
$ws.Range("B34").Value = "Bacillus"
$ws.Range("A35").Value = "bodyCoccus"
$ws.Range("B35").Value = "Coccus"
$ws.Range("A36").Value = "bodyCoccobacillus"
$ws.Range("B36").Value = "Coccobacillus"
$ws.Range("A37").Value = "bodySpirillum"
$ws.Range("B37").Value = "Spirillum"
$ws.Range("A38").Value = "cellStructureThermophile"
$ws.Range("B38").Value = "Thermophile"
$ws.Range("A39").Value = "cellStructurePsychrophile"
$ws.Range("B39").Value = "Psychrophile"
$ws.Range("A40").Value = "cellStructureMethanogen"
$ws.Range("B40").Value = "Methanogen"
$ws.Range("A41").Value = "cellStructureHalophile"
$ws.Range("B41").Value = "Halophile"
$ws.Range("A42").Value = "motilityFlagellaMonotrichous"
$ws.Range("B42").Value = "Monotrichous Flagella"
$ws.Range("A43").Value = "motilityFlagellaLophotrichous"
$ws.Range("B43").Value = "Lophotrichous Flagella"
$ws.Range("A44").Value = "motilityFlagellaPeritrichous"
$ws.Range("B44").Value = "Peritrichous Flagella"
$ws.Range("A45").Value = "motilityFlagellaAmphitrichous"
$ws.Range("B45").Value = "Amphitrichous Flagella"
$ws.Range("A46").Value = "metabolismMethanotroph"
$ws.Range("B46").Value = "Methanotroph"
$ws.Range("A47").Value = "metabolismPhotoautotroph"
$ws.Range("B47").Value = "Photoautotroph"
$ws.Range("A48").Value = "metabolismOrganotroph"
$ws.Range("B48").Value = "Organotroph"
$ws.Range("A49").Value = "metabolismEndobiotic"
$ws.Range("B49").Value = "Organotroph (Endobiotic)"

# --- Insert the new "completed / select / back" strings.

$ws.Range("A17").Value = "completed"
$ws.Range("B17").Value = "COMPLETED"
$ws.Range("A18").Value = "select"
$ws.Range("B18").Value = "SELECT"
$ws.Range("A19").Value = "back"
$ws.Range("B19").Value = "BACK"

# --- categoryBody switches from "Shape" to "Body".

$ws.Range("B27").Value = "Body"

# --- Insert the new "ok" row.

$ws.Range("A16").Value = "ok"
$ws.Range("B16").Value = "OKAY"

# --- Insert the new "retry / changeEnvironment / editOrganism" rows.

$ws.Range("A22").Value = "retry"
$ws.Range("B22").Value = "RETRY"
$ws.Range("A23").Value = "changeEnvironment"
$ws.Range("B23").Value = "CHANGE ENVIRONMENT"
$ws.Range("A24").Value = "editOrganism"
$ws.Range("B24").Value = "EDIT ORGANISM"

# --- Insert the new "timeExpire / timeExpireDesc" rows.

$ws.Range("A20").Value = "timeExpire"
$ws.Range("B20").Value = "Time's Up"
$ws.Range("A21").Value = "timeExpireDesc"
$ws.Range("B21").Value = "You only got {0} out of {1} organisms."

# --- Refresh the sheet view: the visible window now starts scrolled to
#     row 10, with B21 (the new "timeExpireDesc" value cell) active.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B21").Select()
